$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Narrow column B slightly (closest achievable to 38.5703125 character-units) ---
$ws.Columns.Item(2).ColumnWidth = 37.7

# --- Insert 4 new rows before row 20 (new Team-4 issues) ---
$ws.Rows.Item(20).Resize(4).Insert()

# Row 20
$ws.Cells.Item(20,1).Value = 4
$ws.Cells.Item(20,2).Value = 'Inappropriate error message is displayed when changing the root folder (for first time sync on the first computer)'
$ws.Cells.Item(20,3).Value = 'B'
$ws.Cells.Item(20,4).Value = 'Select a root folder by drag & drop. 
Then change the root folder using drag & drop again. A error message "An error has occurred while loading data! Please choose the correct sync folder." is displayed.
This does not happen if I use the "Browse" button to change the root folder'
$ws.Cells.Item(20,5).Value = 'Jiayuan'

# Row 21
$ws.Cells.Item(21,1).Value = 4
$ws.Cells.Item(21,2).Value = 'System reports "We are done synchronizing" but no files are copied to Tzync folder in the exe root folder (For first time sync on the first computer)'
$ws.Cells.Item(21,3).Value = 'B'
$ws.Cells.Item(21,4).Value = 'Select a folder to sync on the first computer.
Change to another folder using "Browse button".
Click on Sync button.
Click OK when confirmation dialog boxes appear.
System reports "We are done synchronizing" but no files are copied to the Tzync folder.'
$ws.Cells.Item(21,5).Value = 'Jiayuan'

# Row 22
$ws.Cells.Item(22,1).Value = 4
$ws.Cells.Item(22,2).Value = 'System displays inappropriate error message when sync is clicked (For first time sync on the first computer)'
$ws.Cells.Item(22,3).Value = 'B'
$ws.Cells.Item(22,4).Value = 'Select a folder to sync on the first computer.
Change to another folder using "Browse button".
Click on Sync button.
Click OK for 2 times.
System display error message "Please choose a valid folder" while the folder input is a valid directory.'
$ws.Cells.Item(22,5).Value = 'Jiayuan'

# Row 23
$ws.Cells.Item(23,1).Value = 4
$ws.Cells.Item(23,2).Value = 'Not able to setup sync for the second computer after sync has been sucessfully performed on the first computer'
$ws.Cells.Item(23,3).Value = 'B'
$ws.Cells.Item(23,4).Value = 'Select a root folder to sync on the first computer.
Click Sync.
Then save the Tzsync folder and the executable on a USB device.
Go to the second computer and run the exe.
Select a folder on the second computer to sync. 
System displays "Ouch that hurts. Don''t do this again" and terminates '
$ws.Cells.Item(23,5).Value = 'Jiayuan'

# --- Apply styles to the 4 new rows (A:center/general, B:left/wrap/text, C:center/wrap/text, D:left/wrap/general, E:left/general) ---
foreach ($r in 20..23) {
  $ws.Cells.Item($r,1).NumberFormat = "General"
  $ws.Cells.Item($r,1).HorizontalAlignment = -4108
  $ws.Cells.Item($r,1).VerticalAlignment = -4108
  $ws.Cells.Item($r,1).WrapText = $false

  $ws.Cells.Item($r,2).NumberFormat = "@"
  $ws.Cells.Item($r,2).HorizontalAlignment = -4131
  $ws.Cells.Item($r,2).VerticalAlignment = -4108
  $ws.Cells.Item($r,2).WrapText = $true

  $ws.Cells.Item($r,3).NumberFormat = "@"
  $ws.Cells.Item($r,3).HorizontalAlignment = -4108
  $ws.Cells.Item($r,3).VerticalAlignment = -4108
  $ws.Cells.Item($r,3).WrapText = $true

  $ws.Cells.Item($r,4).NumberFormat = "General"
  $ws.Cells.Item($r,4).HorizontalAlignment = -4131
  $ws.Cells.Item($r,4).VerticalAlignment = -4108
  $ws.Cells.Item($r,4).WrapText = $true

  $ws.Cells.Item($r,5).NumberFormat = "General"
  $ws.Cells.Item($r,5).HorizontalAlignment = -4131
  $ws.Cells.Item($r,5).VerticalAlignment = -4108
  $ws.Cells.Item($r,5).WrapText = $false
}

# --- Exact row heights for the new rows 20-23 ---
$ws.Rows.Item(20).RowHeight = 105
$ws.Rows.Item(21).RowHeight = 90
$ws.Rows.Item(22).RowHeight = 90
$ws.Rows.Item(23).RowHeight = 120

# --- Append new row 75 (Team 14) at the end ---
$ws.Cells.Item(75,1).Value = 14
$ws.Cells.Item(75,2).Value = 'System displays "Folders that are to be synchronized cannot be a subdirectory of each other" where the source and destination folders are not sub-directory of each other'
$ws.Cells.Item(75,3).Value = 'B'
$ws.Cells.Item(75,4).Value = 'Create a folder named "A".
Create another folder named "A1" with B as its sub-folder.
Select folder "A" as source, then select "A1\B" as destination directory'
$ws.Cells.Item(75,5).Value = 'Jiayuan'

$ws.Cells.Item(75,1).NumberFormat = "General"
$ws.Cells.Item(75,1).HorizontalAlignment = -4108
$ws.Cells.Item(75,1).VerticalAlignment = -4108
$ws.Cells.Item(75,1).WrapText = $false

$ws.Cells.Item(75,2).NumberFormat = "@"
$ws.Cells.Item(75,2).HorizontalAlignment = -4131
$ws.Cells.Item(75,2).VerticalAlignment = -4108
$ws.Cells.Item(75,2).WrapText = $true

$ws.Cells.Item(75,3).NumberFormat = "@"
$ws.Cells.Item(75,3).HorizontalAlignment = -4108
$ws.Cells.Item(75,3).VerticalAlignment = -4108
$ws.Cells.Item(75,3).WrapText = $true

$ws.Cells.Item(75,4).NumberFormat = "General"
$ws.Cells.Item(75,4).HorizontalAlignment = -4131
$ws.Cells.Item(75,4).VerticalAlignment = -4108
$ws.Cells.Item(75,4).WrapText = $true

$ws.Cells.Item(75,5).NumberFormat = "@"
$ws.Cells.Item(75,5).HorizontalAlignment = -4131
$ws.Cells.Item(75,5).VerticalAlignment = -4108
$ws.Cells.Item(75,5).WrapText = $true

$ws.Rows.Item(75).RowHeight = 75

# --- Remove the stale selection left over from before the edit (best effort) ---
$ws.Range("A1").Select()
